# Auto-generated edit script applying numeric value updates to H:N columns
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets, per the commit diff.
$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(12, 8).Value = 446.25
$ws.Cells.Item(12, 9).Value = 392.5
$ws.Cells.Item(12, 10).Value = 457
$ws.Cells.Item(12, 11).Value = 392.5
$ws.Cells.Item(12, 12).Value = 457
$ws.Cells.Item(12, 13).Value = -222.5
$ws.Cells.Item(12, 14).Value = -797
$ws.Cells.Item(18, 8).Value = 17400
$ws.Cells.Item(18, 10).Value = 3000
$ws.Cells.Item(18, 12).Value = 3000
$ws.Cells.Item(18, 14).Value = -3568
$ws.Cells.Item(40, 8).Value = 5393.5483
$ws.Cells.Item(40, 9).Value = 887.5
$ws.Cells.Item(40, 11).Value = 887.5
$ws.Cells.Item(40, 13).Value = -712.5
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(51, 9).Value = 0
$ws.Cells.Item(51, 11).Value = 0
$ws.Cells.Item(51, 13).ClearContents()
$ws.Cells.Item(86, 8).Value = 1000
$ws.Cells.Item(86, 9).Value = 1000
$ws.Cells.Item(86, 11).Value = 1000
$ws.Cells.Item(86, 13).Value = 123
$ws.Cells.Item(89, 8).Value = 1000
$ws.Cells.Item(89, 9).Value = 1000
$ws.Cells.Item(89, 11).Value = 5000
$ws.Cells.Item(89, 13).Value = 616
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()
$ws.Cells.Item(114, 8).Value = 0
$ws.Cells.Item(114, 10).Value = 0
$ws.Cells.Item(114, 12).Value = 0
$ws.Cells.Item(114, 14).ClearContents()
$ws.Cells.Item(130, 8).Value = 97998.75
$ws.Cells.Item(130, 10).Value = 97998.75
$ws.Cells.Item(130, 12).Value = 97998.75
$ws.Cells.Item(130, 14).Value = -108038.75
$ws.Cells.Item(136, 8).Value = 40000
$ws.Cells.Item(136, 10).Value = 40000
$ws.Cells.Item(136, 12).Value = 40000
$ws.Cells.Item(136, 14).Value = -50200
$ws.Cells.Item(139, 8).Value = 99987
$ws.Cells.Item(139, 10).Value = 99987
$ws.Cells.Item(139, 12).Value = 99987
$ws.Cells.Item(139, 14).Value = -110267
$ws.Cells.Item(140, 8).Value = 99995
$ws.Cells.Item(140, 10).Value = 99995
$ws.Cells.Item(140, 12).Value = 99995
$ws.Cells.Item(140, 14).Value = -110355

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(76, 8).Value = 60000
$ws.Cells.Item(76, 10).Value = 60000
$ws.Cells.Item(76, 12).Value = 60000
$ws.Cells.Item(76, 14).Value = -60676
$ws.Cells.Item(79, 8).Value = 60000
$ws.Cells.Item(79, 10).Value = 60000
$ws.Cells.Item(79, 12).Value = 60000
$ws.Cells.Item(79, 14).Value = -62340
$ws.Cells.Item(130, 8).Value = 94000
$ws.Cells.Item(130, 10).Value = 94000
$ws.Cells.Item(130, 12).Value = 94000
$ws.Cells.Item(130, 14).Value = -104040
$ws.Cells.Item(138, 8).Value = 749999.5
$ws.Cells.Item(138, 10).Value = 749999.5
$ws.Cells.Item(138, 12).Value = 749999.5
$ws.Cells.Item(138, 14).Value = -760279.5

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(52, 8).Value = 0
$ws.Cells.Item(52, 10).Value = 0
$ws.Cells.Item(52, 12).Value = 0
$ws.Cells.Item(52, 14).ClearContents()
$ws.Cells.Item(108, 8).Value = 0
$ws.Cells.Item(108, 10).Value = 0
$ws.Cells.Item(108, 12).Value = 0
$ws.Cells.Item(108, 14).ClearContents()
$ws.Cells.Item(120, 8).Value = 99994
$ws.Cells.Item(120, 10).Value = 99994
$ws.Cells.Item(120, 12).Value = 99994
$ws.Cells.Item(120, 14).Value = -109670
$ws.Cells.Item(121, 8).Value = 0
$ws.Cells.Item(121, 10).Value = 0
$ws.Cells.Item(121, 12).Value = 0
$ws.Cells.Item(121, 14).ClearContents()
$ws.Cells.Item(130, 8).Value = 25000
$ws.Cells.Item(130, 10).Value = 25000
$ws.Cells.Item(130, 12).Value = 25000
$ws.Cells.Item(130, 14).Value = -35040
$ws.Cells.Item(137, 8).Value = 55000
$ws.Cells.Item(137, 10).Value = 0
$ws.Cells.Item(137, 12).Value = 0
$ws.Cells.Item(137, 14).ClearContents()
$ws.Cells.Item(140, 8).Value = 99995
$ws.Cells.Item(140, 10).Value = 99995
$ws.Cells.Item(140, 12).Value = 99995
$ws.Cells.Item(140, 14).Value = -110355

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 240.77777
$ws.Cells.Item(7, 10).Value = 288.16666
$ws.Cells.Item(7, 12).Value = 288.16666
$ws.Cells.Item(7, 14).Value = -514.16666
$ws.Cells.Item(22, 8).Value = 250
$ws.Cells.Item(22, 9).Value = 200
$ws.Cells.Item(22, 10).Value = 400
$ws.Cells.Item(22, 11).Value = 200
$ws.Cells.Item(22, 12).Value = 400
$ws.Cells.Item(22, 13).Value = 150
$ws.Cells.Item(22, 14).Value = -1100
$ws.Cells.Item(44, 8).Value = 1982
$ws.Cells.Item(44, 9).Value = 1982
$ws.Cells.Item(44, 11).Value = 1982
$ws.Cells.Item(44, 13).Value = -1540
$ws.Cells.Item(52, 8).Value = 99987
$ws.Cells.Item(52, 10).Value = 99987
$ws.Cells.Item(52, 12).Value = 99987
$ws.Cells.Item(52, 14).Value = -100575
$ws.Cells.Item(100, 8).Value = 52999
$ws.Cells.Item(100, 10).Value = 52999
$ws.Cells.Item(100, 12).Value = 52999
$ws.Cells.Item(100, 14).Value = -55163
$ws.Cells.Item(116, 8).Value = 0
$ws.Cells.Item(116, 10).Value = 0
$ws.Cells.Item(116, 12).Value = 0
$ws.Cells.Item(116, 14).ClearContents()

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(46, 8).Value = 500
$ws.Cells.Item(46, 9).Value = 100
$ws.Cells.Item(46, 10).Value = 633.3333
$ws.Cells.Item(46, 11).Value = 300
$ws.Cells.Item(46, 12).Value = 1899.9999
$ws.Cells.Item(46, 13).Value = -209
$ws.Cells.Item(46, 14).Value = -2081.9999
$ws.Cells.Item(129, 8).Value = 1499.6666
$ws.Cells.Item(129, 9).Value = 999.5
$ws.Cells.Item(129, 10).Value = 2500
$ws.Cells.Item(129, 11).Value = 2998.5
$ws.Cells.Item(129, 12).Value = 7500
$ws.Cells.Item(129, 13).Value = 2001.5
$ws.Cells.Item(129, 14).Value = -17500

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(5, 8).Value = 0
$ws.Cells.Item(5, 9).Value = 0
$ws.Cells.Item(5, 11).Value = 0
$ws.Cells.Item(5, 13).ClearContents()
$ws.Cells.Item(11, 8).Value = 3771428.5
$ws.Cells.Item(11, 9).Value = 3483333.2
$ws.Cells.Item(11, 11).Value = 3483333.2
$ws.Cells.Item(11, 13).Value = -3483194.2
$ws.Cells.Item(122, 8).Value = 3140
$ws.Cells.Item(122, 9).Value = 1400
$ws.Cells.Item(122, 10).Value = 5750
$ws.Cells.Item(122, 11).Value = 4200
$ws.Cells.Item(122, 12).Value = 17250
$ws.Cells.Item(122, 13).Value = -1750
$ws.Cells.Item(122, 14).Value = -22150
$ws.Cells.Item(133, 8).Value = 0
$ws.Cells.Item(133, 10).Value = 0
$ws.Cells.Item(133, 12).Value = 0
$ws.Cells.Item(133, 14).ClearContents()

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(2, 8).Value = 3816.6667
$ws.Cells.Item(2, 10).Value = 3816.6667
$ws.Cells.Item(2, 12).Value = 3816.6667
$ws.Cells.Item(2, 14).Value = -4040.6667
$ws.Cells.Item(3, 8).Value = 12800
$ws.Cells.Item(3, 9).Value = 0
$ws.Cells.Item(3, 11).Value = 0
$ws.Cells.Item(3, 13).ClearContents()
$ws.Cells.Item(15, 8).Value = 12800
$ws.Cells.Item(15, 9).Value = 0
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 13).ClearContents()
$ws.Cells.Item(46, 8).Value = 3237.9092
$ws.Cells.Item(46, 9).Value = 1000.5
$ws.Cells.Item(46, 11).Value = 1000.5
$ws.Cells.Item(46, 13).Value = -812.5
$ws.Cells.Item(93, 8).Value = 6000
$ws.Cells.Item(93, 9).Value = 6000
$ws.Cells.Item(93, 11).Value = 6000
$ws.Cells.Item(93, 13).Value = -4752
$ws.Cells.Item(122, 8).Value = 6600
$ws.Cells.Item(122, 9).Value = 6133.3335
$ws.Cells.Item(122, 11).Value = 18400.0005
$ws.Cells.Item(122, 13).Value = -15950.0005
$ws.Cells.Item(128, 8).Value = 89000
$ws.Cells.Item(128, 10).Value = 89000
$ws.Cells.Item(128, 12).Value = 89000
$ws.Cells.Item(128, 14).Value = -98960
$ws.Cells.Item(135, 8).Value = 30000
$ws.Cells.Item(135, 10).Value = 30000
$ws.Cells.Item(135, 12).Value = 30000
$ws.Cells.Item(135, 14).Value = -40140

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(2, 8).Value = 1569476.1
$ws.Cells.Item(2, 9).Value = 2266799.5
$ws.Cells.Item(2, 10).Value = 498.75
$ws.Cells.Item(2, 11).Value = 2266799.5
$ws.Cells.Item(2, 12).Value = 498.75
$ws.Cells.Item(2, 13).Value = -2266687.5
$ws.Cells.Item(2, 14).Value = -722.75
$ws.Cells.Item(4, 8).Value = 22907.334
$ws.Cells.Item(4, 9).Value = 33714.082
$ws.Cells.Item(4, 10).Value = 1293.8334
$ws.Cells.Item(4, 11).Value = 33714.082
$ws.Cells.Item(4, 12).Value = 1293.8334
$ws.Cells.Item(4, 13).Value = -33601.082
$ws.Cells.Item(4, 14).Value = -1519.8334
$ws.Cells.Item(19, 8).Value = 14502
$ws.Cells.Item(19, 9).Value = 10005
$ws.Cells.Item(19, 11).Value = 10005
$ws.Cells.Item(19, 13).Value = -9831
$ws.Cells.Item(68, 8).Value = 26066.666
$ws.Cells.Item(68, 10).Value = 26066.666
$ws.Cells.Item(68, 12).Value = 26066.666
$ws.Cells.Item(68, 14).Value = -27688.666
$ws.Cells.Item(71, 8).Value = 26066.666
$ws.Cells.Item(71, 10).Value = 26066.666
$ws.Cells.Item(71, 12).Value = 78199.99800000001
$ws.Cells.Item(71, 14).Value = -86311.99800000001
$ws.Cells.Item(112, 8).Value = 50000
$ws.Cells.Item(112, 10).Value = 50000
$ws.Cells.Item(112, 12).Value = 50000
$ws.Cells.Item(112, 14).Value = -52954
$ws.Cells.Item(118, 8).Value = 114988
$ws.Cells.Item(118, 10).Value = 175000
$ws.Cells.Item(118, 12).Value = 175000
$ws.Cells.Item(118, 14).Value = -178314
$ws.Cells.Item(123, 8).Value = 49000
$ws.Cells.Item(123, 10).Value = 49000
$ws.Cells.Item(123, 12).Value = 49000
$ws.Cells.Item(123, 14).Value = -58800
$ws.Cells.Item(135, 8).Value = 99600
$ws.Cells.Item(135, 10).Value = 99600
$ws.Cells.Item(135, 12).Value = 99600
$ws.Cells.Item(135, 14).Value = -109740
